# Append a new log row (row 23) to the Nalco run-log worksheet, matching
# the style of the existing data rows (centered alignment, same cell style
# used throughout the sheet body).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 23

# Match the formatting already used by the rest of the data rows.
$rowRange = $ws.Range("A" + $newRow + ":H" + $newRow)
$rowRange.HorizontalAlignment = -4108   # xlCenter
$rowRange.VerticalAlignment = -4108     # xlCenter

$ws.Cells.Item($newRow, 1).Value = "2025-08-17 06:48:59 UTC"
$ws.Cells.Item($newRow, 2).Value = "2025-08-17 12:18:59 IST"
$ws.Cells.Item($newRow, 3).Value = "SKIPPED"
$ws.Cells.Item($newRow, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($newRow, 5).Value = "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf"
# Column F (Saved PDF) stays blank - no PDF was downloaded for a skipped run.
$ws.Cells.Item($newRow, 7).Value = 0
# Column H (Total Rows After) stays blank, as in the other SKIPPED rows.
